$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy the format from the neighboring header cell G1
# (bold font, border, centered alignment) then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data column values for rows 2-8
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
